$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: remove leading spaces inside quotes and drop the trailing empty ' ' element
$ws.Range("D2").Value = "['beat', 'raid', 'injury', 'rubber', 'stungrenade', 'stuntgrenade', 'tear', 'resisted', 'brutality']"

# F2: collapse the list into a single comma-joined string element, using 'resisted' instead of 'resist'
$ws.Range("F2").Value = "['beat,raid,injury,rubber,stungrenade,stuntgrenade,tear,resisted,brutality']"

# C3: remove leading spaces inside quotes and replace trailing empty '' with 'vigilante'
$ws.Range("C3").Value = "['Beat', 'beating', 'mob justice', 'necklace', 'necklacing', 'vigilantes', 'vigilante']"

# E3: regroup into two comma-joined string elements
$ws.Range("E3").Value = "['beat,beating,mob', 'justice,necklace,necklacing,vigilantes,vigilante']"
